# ReciboDePago.xlsx - "correcciones en la generacion de comprobantes"
#
# The workbook is a single-receipt template ("Hoja1"). This edit fills the
# template with a new receipt's data (replacing the previous sample data):
#   - N°                 (P3)  -> "N°:246"
#   - Sr. (payer)         (L7)  -> "CALIRI PICON DIEGO , DNI 30819184"
#   - Importe             (E8)  -> 1111
#   - Con domicilio en    (L9)  -> "13"            (kept as text)
#   - La cantidad de      (J11) -> "UN MIL CIENTO ONCE PESOS "
#   - Efectivo por/Subtot (P15) -> 1111
#   - Observaciones       (K18) -> "asd1235asdf"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Receipt number
$ws.Range("P3").Value() = "N°:246"

# Payer name / DNI
$ws.Range("L7").Value() = "CALIRI PICON DIEGO , DNI 30819184"

# Amount (numeric)
$ws.Range("E8").Value() = 1111

# Address field holds a numeric-looking value ("13") but must remain text,
# matching the original cell's string type (t="s"). Prefixing with an
# apostrophe forces Excel to store it as text instead of coercing it to a
# number.
$ws.Range("L9").Value() = "'13"

# Amount in words
$ws.Range("J11").Value() = "UN MIL CIENTO ONCE PESOS "

# Subtotal / "Efectivo por" amount (numeric)
$ws.Range("P15").Value() = 1111

# Observaciones
$ws.Range("K18").Value() = "asd1235asdf"
